$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply crypto price / volume(1h) updates scraped on Sun Nov  5 07:45:34 UTC 2023
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '35.537.32'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +1.40%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.903.90'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +3.11%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.64%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '246.52'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +5.83%  '
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +1.61%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '42.12'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +2.99%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0997'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.55%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.180.16'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +3.10%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '12.39'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +8.94%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.932.75'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +4.61%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.690'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +2.06%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.85'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +3.41%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '35.573.53'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +1.47%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0₃0812'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +2.35%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '243.88'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +1.24%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.48'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +2.86%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.88'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +2.28%  '
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.53%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.29'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.61%  '
$ws.Range('B25').NumberFormat = "@"
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').NumberFormat = "@"
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.22'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +28.82%  '
$ws.Range('B26').NumberFormat = "@"
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').NumberFormat = "@"
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '172.13'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.61'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +9.05%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '17.97'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +2.20%  '
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.85%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.983'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +28.62%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.11'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +3.37%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0566'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.63%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.17'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +4.92%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.74'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +5.89%  '
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +2.56%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.33'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +7.25%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.10'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +2.34%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +2.33%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '91.03'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -0.42%  '
$ws.Range('B41').NumberFormat = "@"
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').NumberFormat = "@"
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.352.91'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.49%  '
$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '15.51'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +5.84%  '
$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '50.07'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +45.21%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0592'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +11.10%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.36'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +1.78%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.77'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.41'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +0.96%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '6.66'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +4.88%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.75'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.092.27'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +3.10%  '
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +2.43%  '
